# Updates the cryptocurrency price (column D) and 1h volume-change
# (column E) values on the active worksheet to reflect the latest scrape,
# per the commit "Updated cryptos list ... with GitHub Actions".
#
# Column D values are stored as plain text in the workbook (e.g. "1.01",
# "310.08") rather than numbers, matching the scraper's original inline
# string output -- so for any new price that Excel would otherwise
# auto-detect as a number, the cell is pre-formatted as Text ("@") before
# the value is written, keeping it a text string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.677.37"
$ws.Cells.Item(2, 5).Value = "  +0.53%  "
$ws.Cells.Item(3, 4).Value = "2.300.67"
$ws.Cells.Item(3, 5).Value = "  +0.01%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.01"
$ws.Cells.Item(4, 5).Value = "  +0.89%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "310.08"
$ws.Cells.Item(5, 5).Value = "  -1.85%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "105.43"
$ws.Cells.Item(6, 5).Value = "  +1.59%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.627"
$ws.Cells.Item(7, 5).Value = "  -0.22%  "
$ws.Cells.Item(8, 5).Value = "  +0.08%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.605"
$ws.Cells.Item(9, 5).Value = "  -0.06%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "39.60"
$ws.Cells.Item(10, 5).Value = "  -0.49%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0907"
$ws.Cells.Item(11, 5).Value = "  +0.25%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.24"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.107"
$ws.Cells.Item(13, 5).Value = "  +0.07%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.989"
$ws.Cells.Item(14, 5).Value = "  -0.69%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.24"
$ws.Cells.Item(15, 5).Value = "  -0.70%  "
$ws.Cells.Item(16, 4).Value = "2.652.86"
$ws.Cells.Item(16, 5).Value = "  +0.09%  "
$ws.Cells.Item(17, 4).Value = "2.300.26"
$ws.Cells.Item(17, 5).Value = "  -1.10%  "
$ws.Cells.Item(18, 4).Value = "42.858.81"
$ws.Cells.Item(18, 5).Value = "  +0.70%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.31"
$ws.Cells.Item(19, 5).Value = "  -4.01%  "
$ws.Cells.Item(20, 5).Value = "  -1.01%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "13.42"
$ws.Cells.Item(21, 5).Value = "  -1.31%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "73.57"
$ws.Cells.Item(22, 5).Value = "  -0.59%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.46"
$ws.Cells.Item(23, 5).Value = "  -2.32%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "266.91"
$ws.Cells.Item(24, 5).Value = "  -0.33%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.23"
$ws.Cells.Item(25, 5).Value = "  -0.21%  "
$ws.Cells.Item(26, 5).Value = "  +0.07%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.50"
$ws.Cells.Item(27, 5).Value = "  +12.38%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "10.95"
$ws.Cells.Item(28, 5).Value = "  +0.54%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.29"
$ws.Cells.Item(29, 5).Value = "  -2.57%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "37.92"
$ws.Cells.Item(30, 5).Value = "  +0.66%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "22.23"
$ws.Cells.Item(31, 5).Value = "  -1.46%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "164.88"
$ws.Cells.Item(32, 5).Value = "  -0.37%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0865"
$ws.Cells.Item(33, 5).Value = "  -2.03%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.84"
$ws.Cells.Item(34, 5).Value = "  +6.53%  "
$ws.Cells.Item(35, 5).Value = "  -0.58%  "
$ws.Cells.Item(36, 5).Value = "  -1.10%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.62"
$ws.Cells.Item(37, 5).Value = "  +0.75%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0353"
$ws.Cells.Item(38, 5).Value = "  -0.26%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.81"
$ws.Cells.Item(39, 5).Value = "  +2.76%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.62"
$ws.Cells.Item(40, 5).Value = "  -3.35%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "107.10"
$ws.Cells.Item(41, 5).Value = "  +9.02%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.56"
$ws.Cells.Item(42, 5).Value = "  -3.18%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "71.37"
$ws.Cells.Item(43, 5).Value = "  +1.90%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.228"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.01"
$ws.Cells.Item(45, 5).Value = "  +0.20%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "12.30"
$ws.Cells.Item(46, 5).Value = "  -0.49%  "
$ws.Cells.Item(47, 4).Value = "1.709.19"
$ws.Cells.Item(47, 5).Value = "  +4.42%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "111.48"
$ws.Cells.Item(48, 5).Value = "  -4.46%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "76.28"
$ws.Cells.Item(49, 5).Value = "  -5.55%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "8.87"
$ws.Cells.Item(50, 5).Value = "  -0.71%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "5.19"
$ws.Cells.Item(51, 5).Value = "  -1.91%  "
